# repull data, push all data, mean calculation
# Updates column F (dSF) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -9
    8  = -4
    9  = -4
    10 = -3
    11 = 0
    14 = -2
    15 = -3
    18 = -3
    19 = 0
    23 = 11
    24 = -4
    26 = -3
    27 = -1
    29 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
